$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "37.840.39"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.18%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.027.26"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  +0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "227.35"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("E6").Value = "  -0.28%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "59.22"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +1.76%  "

$ws.Range("E8").Value = "  -0.01%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.383"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "2.330.25"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "14.53"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -0.27%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "21.06"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +1.67%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.758"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("E16").Value = "  -2.11%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.029.88"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -0.52%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "37.741.81"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.28%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.01"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -2.08%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "69.95"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +0.43%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0821"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -1.20%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "224.95"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  -1.91%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.20"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -2.02%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.27"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.11%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "165.15"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -0.76%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.128"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -3.31%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "18.91"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("E30").Value = "  -4.84%  "

$ws.Range("E31").Value = "  +0.72%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.42"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -2.50%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +1.23%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.50"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -1.61%  "

$ws.Range("E35").Value = "  -1.62%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.36"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +6.62%  "

$ws.Range("E37").Value = "  -3.16%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("E39").Value = "  +0.13%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.519.54"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.30%  "

$ws.Range("E41").Value = "  +0.29%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "96.41"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -1.35%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "16.73"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +0.56%  "

$ws.Range("E44").Value = "  -0.59%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0915"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("E46").Value = "  -1.83%  "

$ws.Range("E47").Value = "  -4.08%  "

$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("E49").Value = "  -0.40%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "7.05"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +0.75%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.218.43"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.16%  "
